$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.872.69"
$ws.Range('E2').Value = '  +6.25%  '
$ws.Range("D3").Value = "'1.731.87"
$ws.Range('E3').Value = '  +4.59%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range("D5").Value = "'229.28"
$ws.Range("D6").Value = "'0.5428"
$ws.Range('E6').Value = '  +3.76%  '
$ws.Range("D7").Value = "'1.002"
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range("D8").Value = "'0.2771"
$ws.Range('E8').Value = '  +3.99%  '
$ws.Range("D9").Value = "'0.06774"
$ws.Range('E9').Value = '  +6.50%  '
$ws.Range("D10").Value = "'21.67"
$ws.Range('E10').Value = '  +5.15%  '
$ws.Range("D11").Value = "'0.07828"
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range("D12").Value = "'4.714"
$ws.Range('E12').Value = '  +2.65%  '
$ws.Range("D13").Value = "'1.779.66"
$ws.Range('E13').Value = '  +8.66%  '
$ws.Range("D14").Value = "'1.970.16"
$ws.Range('E14').Value = '  +4.58%  '
$ws.Range("D15").Value = "'0.6017"
$ws.Range('E15').Value = '  +6.60%  '
$ws.Range("D16").Value = "'0.0₅8394"
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range("D17").Value = "'68.76"
$ws.Range('E17').Value = '  +5.13%  '
$ws.Range("D18").Value = "'27.845.41"
$ws.Range('E18').Value = '  +6.17%  '
$ws.Range("D19").Value = "'215.93"
$ws.Range('E19').Value = '  +12.56%  '
$ws.Range("D20").Value = "'4.827"
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range("D21").Value = "'1.000"
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('E22').Value = '  +5.23%  '
$ws.Range("D23").Value = "'6.249"
$ws.Range('E23').Value = '  +3.97%  '
$ws.Range("D24").Value = "'1.002"
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range("D25").Value = "'145.53"
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range("D26").Value = "'0.1249"
$ws.Range('E26').Value = '  +4.08%  '
$ws.Range("D27").Value = "'7.441"
$ws.Range('E27').Value = '  +2.24%  '
$ws.Range("D28").Value = "'1.641"
$ws.Range('E28').Value = '  +9.23%  '
$ws.Range('E29').Value = '  +5.77%  '
$ws.Range("D30").Value = "'0.05606"
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range("D31").Value = "'1.318"
$ws.Range('E31').Value = '  +3.26%  '
$ws.Range("D32").Value = "'3.735"
$ws.Range('E32').Value = '  +6.72%  '
$ws.Range("D33").Value = "'3.531"
$ws.Range('E33').Value = '  +5.29%  '
$ws.Range("D34").Value = "'1.636"
$ws.Range('E34').Value = '  +3.60%  '
$ws.Range("D35").Value = "'0.9821"
$ws.Range('E35').Value = '  +3.86%  '
$ws.Range("D36").Value = "'2.856"
$ws.Range('E36').Value = '  +1.88%  '
$ws.Range("D37").Value = "'2.443"
$ws.Range('E37').Value = '  +1.31%  '
$ws.Range("D38").Value = "'0.5939"
$ws.Range('E38').Value = '  +3.31%  '
$ws.Range('E39').Value = '  +4.24%  '
$ws.Range("D40").Value = "'5.949"
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'0.8433"
$ws.Range('E41').Value = '  -0.40%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = "'1.040.83"
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range("D43").Value = "'1.001"
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range("D44").Value = "'102.84"
$ws.Range('E44').Value = '  +1.54%  '
$ws.Range("D45").Value = "'1.874.76"
$ws.Range('E45').Value = '  +4.44%  '
$ws.Range("D46").Value = "'59.82"
$ws.Range('E46').Value = '  +2.47%  '
$ws.Range("D47").Value = "'0.0₈108"
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range("D48").Value = "'8.250"
$ws.Range('E48').Value = '  +2.80%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = "'0.4419"
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.05313"
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('B51').Value = 'Frax'
$ws.Range('C51').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D51").Value = "'0.9981"
$ws.Range('E51').Value = '  -0.49%  '
